# Edit: Fruta / hortaliza, semanal
# The weekly refresh prepends a new observation (row 15) and pushes every
# existing observation down by one row, dropping the previously last
# row's data into the newly appended row 34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Create the brand new row 34 (was not present before) -------------
# Copy the constant / descriptive columns from row 33 (same for every row
# in this sheet) and give D34 the same date/number style as the other
# date cells (numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("A34").Value = 8
$ws.Range("B34").Value = "Terminal La Palmera de La Serena"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100112052
$ws.Range("G34").Value = "Albahaca"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("N34").Value = "`$/paquete"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"

# --- 2. Shift the weekly price/date data down by one row (15 -> 34) ------
# Row 15 receives the freshly reported week; every other row receives the
# values that used to belong to the row above it.

$ws.Range("D15").Value = 44452
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 4500
$ws.Range("L15").Value = 5000
$ws.Range("M15").Value = 4750
$ws.Range("P15").Value = 4750

$ws.Range("D16").Value = 44397
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 4500
$ws.Range("M16").Value = 4250
$ws.Range("P16").Value = 4250

$ws.Range("D17").Value = 44335
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3500
$ws.Range("M17").Value = 3250
$ws.Range("P17").Value = 3250

$ws.Range("D18").Value = 44434
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 4500
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 4750
$ws.Range("P18").Value = 4750

$ws.Range("D19").Value = 44420
$ws.Range("J19").Value = 900
$ws.Range("K19").Value = 4500
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 4750
$ws.Range("P19").Value = 4750

$ws.Range("D20").Value = 44348
$ws.Range("J20").Value = 700
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = 3250
$ws.Range("P20").Value = 3250

$ws.Range("D21").Value = 44427
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = 4500
$ws.Range("L21").Value = 5000
$ws.Range("M21").Value = 4750
$ws.Range("P21").Value = 4750

$ws.Range("D22").Value = 44341
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 3500
$ws.Range("M22").Value = 3250
$ws.Range("P22").Value = 3250

$ws.Range("D23").Value = 44176
$ws.Range("J23").Value = 2000
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 3500
$ws.Range("M23").Value = 3250
$ws.Range("P23").Value = 3250

$ws.Range("D24").Value = 44441
$ws.Range("J24").Value = 600
$ws.Range("K24").Value = 4500
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = 4750
$ws.Range("P24").Value = 4750

$ws.Range("D25").Value = 44432
$ws.Range("J25").Value = 900
$ws.Range("K25").Value = 4500
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = 4750
$ws.Range("P25").Value = 4750

$ws.Range("D26").Value = 44342
$ws.Range("J26").Value = 560
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 3500
$ws.Range("M26").Value = 3250
$ws.Range("P26").Value = 3250

$ws.Range("D27").Value = 44379
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 4000
$ws.Range("L27").Value = 4500
$ws.Range("M27").Value = 4250
$ws.Range("P27").Value = 4250

$ws.Range("D28").Value = 44315
$ws.Range("J28").Value = 700
$ws.Range("K28").Value = 2500
$ws.Range("L28").Value = 3000
$ws.Range("M28").Value = 2750
$ws.Range("P28").Value = 2750

$ws.Range("D29").Value = 44446
$ws.Range("J29").Value = 800
$ws.Range("K29").Value = 4500
$ws.Range("L29").Value = 5000
$ws.Range("M29").Value = 4750
$ws.Range("P29").Value = 4750

$ws.Range("D30").Value = 44411
$ws.Range("J30").Value = 880
$ws.Range("K30").Value = 4000
$ws.Range("L30").Value = 4500
$ws.Range("M30").Value = 4250
$ws.Range("P30").Value = 4250

$ws.Range("D31").Value = 44449
$ws.Range("J31").Value = 700
$ws.Range("K31").Value = 4000
$ws.Range("L31").Value = 4500
$ws.Range("M31").Value = 4250
$ws.Range("P31").Value = 4250

$ws.Range("D32").Value = 44435
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 4500
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = 4750
$ws.Range("P32").Value = 4750

$ws.Range("D33").Value = 44314
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 2500
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = 2750
$ws.Range("P33").Value = 2750

$ws.Range("D34").Value = 44448
$ws.Range("J34").Value = 640
$ws.Range("K34").Value = 4500
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = 4750
$ws.Range("P34").Value = 4750
